# Applies the "Updated cryptos list" commit: refreshes the Price (D) and
# Volume(1h) (E) columns for each coin row, and for rows 42/43 swaps
# MultiversX <-> Algorand (both the coin metadata and their latest figures).
#
# Columns D/E are stored as plain text in the sheet (e.g. "247.64",
# "  -1.03%  "), not numbers, so price-looking values like "247.64" would
# otherwise be auto-coerced to a numeric cell by Excel's Value setter. For
# those cells we briefly force a Text number format, assign the literal
# string, then restore the Normal style so no extra formatting is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.338.22"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.252.18"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0940"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("D14").Value = "2.588.67"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.855"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").Value = "2.254.87"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "42.189.17"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Value = "0.0₃0980"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +26.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("E32").Value = "  -6.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("E34").Value = "  -5.60%  "
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.26%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.203"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  -4.20%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  -2.69%  "
